$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1836.1111
$ws.Range("I28").Value = 1830.4
$ws.Range("J28").Value = 1864.6666
$ws.Range("K28").Value = 1830.4
$ws.Range("L28").Value = 1864.6666
$ws.Range("M28").Value = -1345.4
$ws.Range("N28").Value = -2834.6666
# Row 69
$ws.Range("H69").Value = 13999.667
$ws.Range("J69").Value = 19999.5
$ws.Range("L69").Value = 59998.5
$ws.Range("N69").Value = -61746.5
# Row 72
$ws.Range("H72").Value = 13999.667
$ws.Range("J72").Value = 19999.5
$ws.Range("L72").Value = 179995.5
$ws.Range("N72").Value = -188731.5
# Row 74
$ws.Range("H74").Value = 185365.5
$ws.Range("I74").Value = 503748.5
$ws.Range("J74").Value = 26174
$ws.Range("K74").Value = 503748.5
$ws.Range("L74").Value = 26174
$ws.Range("M74").Value = -502812.5
$ws.Range("N74").Value = -28046
# Row 77
$ws.Range("H77").Value = 185365.5
$ws.Range("I77").Value = 503748.5
$ws.Range("J77").Value = 26174
$ws.Range("K77").Value = 2518742.5
$ws.Range("L77").Value = 130870
$ws.Range("M77").Value = -2514062.5
$ws.Range("N77").Value = -140230
# Row 100
$ws.Range("H100").Value = 1350.5555
$ws.Range("I100").Value = 1458.1666
$ws.Range("J100").Value = 1135.3334
$ws.Range("K100").Value = 1458.1666
$ws.Range("L100").Value = 1135.3334
$ws.Range("M100").Value = -917.1666
$ws.Range("N100").Value = -2217.3334
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 132
$ws.Range("H132").Value = 1484.4419
$ws.Range("I132").Value = 1198.1351
$ws.Range("K132").Value = 3594.4053
$ws.Range("M132").Value = -1064.4053

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2859.4568
$ws.Range("I32").Value = 2568.1184
$ws.Range("J32").Value = 7287.8
$ws.Range("K32").Value = 2568.1184
$ws.Range("L32").Value = 7287.8
$ws.Range("M32").Value = -2281.1184
$ws.Range("N32").Value = -7861.8
# Row 55
$ws.Range("H55").Value = 11249.571
$ws.Range("I55").Value = 7749.8
$ws.Range("K55").Value = 7749.8
$ws.Range("M55").Value = -7434.8
# Row 61
$ws.Range("H61").Value = 3431.3416
$ws.Range("I61").Value = 2718.4285
$ws.Range("J61").Value = 4966.846
$ws.Range("K61").Value = 2718.4285
$ws.Range("L61").Value = 4966.846
$ws.Range("M61").Value = -2506.4285
$ws.Range("N61").Value = -5390.846
# Row 74
$ws.Range("H74").Value = 16898.5
$ws.Range("I74").Value = 1783.8572
$ws.Range("K74").Value = 1783.8572
$ws.Range("M74").Value = -909.8571999999999
# Row 77
$ws.Range("H77").Value = 16898.5
$ws.Range("I77").Value = 1783.8572
$ws.Range("K77").Value = 8919.286
$ws.Range("M77").Value = -4551.286
# Row 122
$ws.Range("H122").Value = 2228.875
$ws.Range("I122").Value = 2266.6191
$ws.Range("K122").Value = 6799.8573
$ws.Range("M122").Value = -4349.8573
# Row 136
$ws.Range("H136").Value = 3431.3416
$ws.Range("I136").Value = 2718.4285
$ws.Range("J136").Value = 4966.846
$ws.Range("K136").Value = 8155.2855
$ws.Range("L136").Value = 14900.538
$ws.Range("M136").Value = -5605.2855
$ws.Range("N136").Value = -20000.538

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 34279.5
$ws.Range("I82").Value = 30024.4
$ws.Range("J82").Value = 55555
$ws.Range("K82").Value = 30024.4
$ws.Range("L82").Value = 55555
$ws.Range("M82").Value = -29641.4
$ws.Range("N82").Value = -56321
# Row 85
$ws.Range("H85").Value = 34279.5
$ws.Range("I85").Value = 30024.4
$ws.Range("J85").Value = 55555
$ws.Range("K85").Value = 30024.4
$ws.Range("L85").Value = 55555
$ws.Range("M85").Value = -28698.4
$ws.Range("N85").Value = -58207
# Row 134
$ws.Range("H134").Value = 1814.9149
$ws.Range("I134").Value = 1779.591
$ws.Range("J134").Value = 2333
$ws.Range("K134").Value = 5338.772999999999
$ws.Range("L134").Value = 6999
$ws.Range("M134").Value = -2803.772999999999
$ws.Range("N134").Value = -12069

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 21090.697
$ws.Range("I31").Value = 27593.205
$ws.Range("J31").Value = 2976.5715
$ws.Range("K31").Value = 27593.205
$ws.Range("L31").Value = 2976.5715
$ws.Range("M31").Value = -27298.205
$ws.Range("N31").Value = -3566.5715
# Row 34
$ws.Range("H34").Value = 21090.697
$ws.Range("I34").Value = 27593.205
$ws.Range("J34").Value = 2976.5715
$ws.Range("K34").Value = 27593.205
$ws.Range("L34").Value = 2976.5715
$ws.Range("M34").Value = -27391.205
$ws.Range("N34").Value = -3380.5715
# Row 59
$ws.Range("H59").Value = 23583.084
# Row 82
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 10000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -9639
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 10000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -8752
$ws.Range("N85").ClearContents()
# Row 134
$ws.Range("H134").Value = 27430.617
$ws.Range("I134").Value = 24060.852
$ws.Range("J134").Value = 40428.285
$ws.Range("K134").Value = 72182.556
$ws.Range("L134").Value = 121284.855
$ws.Range("M134").Value = -69647.556
$ws.Range("N134").Value = -126354.855

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 687.2
$ws.Range("I23").Value = 498
$ws.Range("J23").Value = 734.5
$ws.Range("K23").Value = 1494
$ws.Range("L23").Value = 2203.5
$ws.Range("M23").Value = -1259
$ws.Range("N23").Value = -2673.5
# Row 34
$ws.Range("H34").Value = 2432948.2
$ws.Range("I34").Value = 842088.3
$ws.Range("J34").Value = 3705636.2
$ws.Range("K34").Value = 2526264.9
$ws.Range("L34").Value = 11116908.6
$ws.Range("M34").Value = -2526180.9
$ws.Range("N34").Value = -11117076.6
# Row 55
$ws.Range("H55").Value = 1508.75
# Row 75
$ws.Range("H75").Value = 987
$ws.Range("I75").Value = 987
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2961
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -1963
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 987
$ws.Range("I78").Value = 987
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 8883
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -3891
$ws.Range("N78").ClearContents()
# Row 117
$ws.Range("H117").Value = 665
$ws.Range("J117").Value = 497.5
$ws.Range("L117").Value = 1492.5
$ws.Range("N117").Value = -8376.5
# Row 137
$ws.Range("H137").Value = 4382.222
$ws.Range("J137").Value = 5876.6
$ws.Range("L137").Value = 17629.8
$ws.Range("N137").Value = -27829.8

$ws = $wb.Worksheets.Item("GSM")
# Row 74
$ws.Range("H74").Value = 33332.668
$ws.Range("J74").Value = 33332.668
$ws.Range("L74").Value = 33332.668
$ws.Range("N74").Value = -35204.668
# Row 77
$ws.Range("H77").Value = 33332.668
$ws.Range("J77").Value = 33332.668
$ws.Range("L77").Value = 99998.00399999999
$ws.Range("N77").Value = -109358.004
# Row 127
$ws.Range("H127").Value = 75997.2
$ws.Range("J127").Value = 75997.2
$ws.Range("L127").Value = 75997.2
$ws.Range("N127").Value = -85917.2

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2820.9167
$ws.Range("I61").Value = 2901.8333
$ws.Range("K61").Value = 2901.8333
$ws.Range("M61").Value = -2699.8333
# Row 76
$ws.Range("H76").Value = 14955.25
$ws.Range("J76").Value = 15273.667
$ws.Range("L76").Value = 15273.667
$ws.Range("N76").Value = -15949.667
# Row 79
$ws.Range("H79").Value = 14955.25
$ws.Range("J79").Value = 15273.667
$ws.Range("L79").Value = 15273.667
$ws.Range("N79").Value = -17613.667
# Row 113
$ws.Range("H113").Value = 2820.9167
$ws.Range("I113").Value = 2901.8333
$ws.Range("K113").Value = 2901.8333
$ws.Range("M113").Value = -731.8332999999998

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 6953.5
$ws.Range("I62").Value = 5430.25
$ws.Range("K62").Value = 5430.25
$ws.Range("M62").Value = -4806.25
# Row 65
$ws.Range("H65").Value = 6953.5
$ws.Range("I65").Value = 5430.25
$ws.Range("K65").Value = 27151.25
$ws.Range("M65").Value = -24031.25
# Row 81
$ws.Range("H81").Value = 6493.207
$ws.Range("I81").Value = 11502
$ws.Range("J81").Value = 3857
$ws.Range("K81").Value = 23004
$ws.Range("L81").Value = 7714
$ws.Range("M81").Value = -21943
$ws.Range("N81").Value = -9836
# Row 84
$ws.Range("H84").Value = 6493.207
$ws.Range("I84").Value = 11502
$ws.Range("J84").Value = 3857
$ws.Range("K84").Value = 115020
$ws.Range("L84").Value = 38570
$ws.Range("M84").Value = -109716
$ws.Range("N84").Value = -49178
# Row 113
$ws.Range("H113").Value = 2766.5
$ws.Range("I113").Value = 2266.6667
$ws.Range("J113").Value = 3266.3333
$ws.Range("K113").Value = 6800.000100000001
$ws.Range("L113").Value = 9798.999899999999
$ws.Range("M113").Value = -4630.000100000001
$ws.Range("N113").Value = -14138.9999
# Row 136
$ws.Range("H136").Value = 1684.5312
$ws.Range("I136").Value = 1423.2693
$ws.Range("J136").Value = 2816.6667
$ws.Range("K136").Value = 4269.8079
$ws.Range("L136").Value = 8450.000100000001
$ws.Range("M136").Value = -1719.8079
$ws.Range("N136").Value = -13550.0001
# Row 137
$ws.Range("H137").Value = 98987
$ws.Range("J137").Value = 98987
$ws.Range("L137").Value = 98987
$ws.Range("N137").Value = -109187
